$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35
$ws.Range("A35").Value = 111675581
$ws.Range("I35").Value = "'1"
$ws.Range("I35").Style = "Normal"
$ws.Range("Q35").Value = 690413.7262835158
$ws.Range("R35").Value = 6661427.29424896

# Row 36
$ws.Range("A36").Value = 111675575
$ws.Range("I36").Value = "'1"
$ws.Range("I36").Style = "Normal"
$ws.Range("Q36").Value = 690480.7418955797
$ws.Range("R36").Value = 6661091.463633558

# Row 37
$ws.Range("A37").Value = 111675573
$ws.Range("Q37").Value = 690487.9917822112
$ws.Range("R37").Value = 6661106.352564453

# Row 38
$ws.Range("A38").Value = 111675579
$ws.Range("I38").Value = "'3"
$ws.Range("I38").Style = "Normal"
$ws.Range("Q38").Value = 690425.8424831247
$ws.Range("R38").Value = 6661357.862056008

# Row 39
$ws.Range("A39").Value = 111675580
$ws.Range("I39").Value = "'3"
$ws.Range("I39").Style = "Normal"
$ws.Range("Q39").Value = 690370.5537696742
$ws.Range("R39").Value = 6661292.946251329

# Row 40
$ws.Range("A40").Value = 111675577
$ws.Range("Q40").Value = 690430.9193086301
$ws.Range("R40").Value = 6661356.623615522

# Row 41
$ws.Range("A41").Value = 111675572
$ws.Range("I41").Value = "'2"
$ws.Range("I41").Style = "Normal"
$ws.Range("Q41").Value = 690494.5947179901
$ws.Range("R41").Value = 6661104.692649405

# Row 42
$ws.Range("A42").Value = 111675583
$ws.Range("Q42").Value = 690415.8809986882
$ws.Range("R42").Value = 6661424.403280765

# Row 43
$ws.Range("A43").Value = 111675584
$ws.Range("I43").Value = "'2"
$ws.Range("I43").Style = "Normal"
$ws.Range("Q43").Value = 690414.984509701
$ws.Range("R43").Value = 6661422.355185229

# Row 44
$ws.Range("A44").Value = 111675571
$ws.Range("I44").Value = "'5"
$ws.Range("I44").Style = "Normal"
$ws.Range("Q44").Value = 690509.4285896254
$ws.Range("R44").Value = 6661040.900344189

# Row 45
$ws.Range("A45").Value = 111675582
$ws.Range("Q45").Value = 690352.3333891984
$ws.Range("R45").Value = 6661470.655078794

# Row 46
$ws.Range("A46").Value = 111675578
$ws.Range("Q46").Value = 690368.3990222017
$ws.Range("R46").Value = 6661295.837351476

# Row 47
$ws.Range("A47").Value = 111675574
$ws.Range("Q47").Value = 690486.6986671695
$ws.Range("R47").Value = 6661102.281881573
